# Completed the Add course button for Admin
# Adds the new course row (Maths) to the "courses" sheet and refreshes the
# "No. of courses" counter, then leaves the workbook focused on that sheet.

$wb = $excel.ActiveWorkbook

# settle the previous sheet's selection to a single cell before navigating away
$wsPrev = $wb.Worksheets.Item("removed_students")
$wsPrev.Activate() | Out-Null
$wsPrev.Range("I12").Select() | Out-Null

$ws = $wb.Worksheets.Item("courses")
$ws.Activate() | Out-Null

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Maths"
$ws.Range("C2").Value = "5 Days"
$ws.Range("D2").Value = "Basic Algebra"
$ws.Range("E2").Value = "Sam Davis"

# refresh the course count
$ws.Range("F4").Value = 1

$ws.Range("F4").Select() | Out-Null
